$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.28101533333333
$ws.Range("N2").Value = 36.843046
$ws.Range("O2").Value = 0.959552102275422
$ws.Range("P2").Value = 0.959552102275422
$ws.Range("Q2").Value = 0.2891646933668889
$ws.Range("R2").Value = 2.602482240302
$ws.Range("S2").Value = 0.002701672095932764
$ws.Range("T2").Value = 0.002701672095932764
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("M3").Value = 0.09168666666666665
$ws.Range("O3").Value = 0.007163750827004844
$ws.Range("P3").Value = 0.007163750827004845
$ws.Range("Q3").Value = 0.002158823691111111
$ws.Range("R3").Value = 0.01942941322
$ws.Range("S3").Value = 0.00002016993727139895
$ws.Range("T3").Value = 0.00002016993727139895
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("M4").Value = 0.4259936666666666
$ws.Range("O4").Value = 0.03328414689757318
$ws.Range("P4").Value = 0.03328414689757318
$ws.Range("Q4").Value = 0.01003030487744444
$ws.Range("R4").Value = 0.09027274389699999
$ws.Range("S4").Value = 0.00009371335928175562
$ws.Range("T4").Value = 0.00009371335928175561
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.28101533333333
$ws.Range("N5").Value = 36.843046
$ws.Range("O5").Value = 0.959552102275422
$ws.Range("P5").Value = 0.959552102275422
$ws.Range("Q5").Value = 101.3515557097589
$ws.Range("R5").Value = 912.16400138783
$ws.Range("S5").Value = 0.9469298161965181
$ws.Range("T5").Value = 0.9469298161965179
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("M6").Value = 0.09168666666666665
$ws.Range("O6").Value = 0.007163750827004844
$ws.Range("P6").Value = 0.007163750827004845
$ws.Range("Q6").Value = 0.756662706811111
$ws.Range("S6").Value = 0.007069516327260623
$ws.Range("T6").Value = 0.007069516327260623
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("M7").Value = 0.4259936666666666
$ws.Range("O7").Value = 0.03328414689757318
$ws.Range("P7").Value = 0.03328414689757318
$ws.Range("Q7").Value = 3.515598642889444
$ws.Range("S7").Value = 0.03284631551453813
$ws.Range("T7").Value = 0.03284631551453813
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.28101533333333
$ws.Range("N8").Value = 36.843046
$ws.Range("O8").Value = 0.959552102275422
$ws.Range("P8").Value = 0.959552102275422
$ws.Range("Q8").Value = 1.061820679391778
$ws.Range("R8").Value = 9.556386114525999
$ws.Range("S8").Value = 0.009920613982971192
$ws.Range("T8").Value = 0.009920613982971192
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("M9").Value = 0.09168666666666665
$ws.Range("O9").Value = 0.007163750827004844
$ws.Range("P9").Value = 0.007163750827004845
$ws.Range("Q9").Value = 0.00792725976222222
$ws.Range("R9").Value = 0.07134533785999998
$ws.Range("S9").Value = 0.00007406456247282202
$ws.Range("T9").Value = 0.00007406456247282204
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("M10").Value = 0.4259936666666666
$ws.Range("O10").Value = 0.03328414689757318
$ws.Range("P10").Value = 0.03328414689757318
$ws.Range("Q10").Value = 0.03683155441788888
$ws.Range("R10").Value = 0.3314839897609999
$ws.Range("S10").Value = 0.0003441180237532886
$ws.Range("T10").Value = 0.0003441180237532886
